$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Falta de conhecimento e tecnologia adequado..." paragraph.
# The run containing " e tecnologia" becomes just " " (a single space) and
# the "_GoBack" bookmark (Word drops one at the last edit point) is left
# right after it.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" e tecnologia", $true, $false, $false, $false, $false, `
                   $true, 1, $false, " ", 2)

# Re-anchor the "_GoBack" bookmark to sit right after the text we just
# touched (adding a bookmark with a name that already exists simply moves
# it, which mirrors what Word does when you edit text at that spot).
$goBackRange = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# ---------------------------------------------------------------------------
# Change 3: collapse the three runs that spell out the "vasos auto
# irrigáveis" sentence into a single contiguous run (no text change, just a
# formatting/run no-op edit so the identically-formatted runs coalesce).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("auto irrigáveis", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "auto irrigáveis", 2)
